$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (prices and 1h volume change %) based on refreshed source data.
# Price cells are text-formatted (not numbers) in the original sheet, so we force
# NumberFormat "@" before assigning numeric-looking strings to avoid Excel auto-converting
# them into floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.945.34"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.011.27"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.52"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.606"
$ws.Range("E6").Value = "  -2.03%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.45"
$ws.Range("E8").Value = "  -4.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.378"
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0781"
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("E11").Value = "  -3.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.309.67"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.24"
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.27"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.12"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.007.53"
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.862.33"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.18"
$ws.Range("E19").Value = "  +4.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.70"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.57"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("E25").Value = "  -5.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.96"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.125"
$ws.Range("E28").Value = "  -2.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.65"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("E31").Value = "  -3.13%  "
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0613"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.40"
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("E35").Value = "  -3.78%  "
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.14"
$ws.Range("E38").Value = "  -3.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.35"
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("E40").Value = "  -3.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.477.56"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.77"
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.07"
$ws.Range("E43").Value = "  -2.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0917"
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.77"
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("E46").Value = "  -3.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.17"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.199.66"
$ws.Range("E50").Value = "  -2.02%  "

# Row 51: MultiversX replaced by FTXToken
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.56"
$ws.Range("E51").Value = "  -11.04%  "
